$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.790.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.596.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.14%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.60%  "
$ws.Range("E6").Value = "  -1.86%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -1.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0618"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0839"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.820.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.593.39"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.753.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "209.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.38%  "
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("E21").Value = "  -1.60%  "
$ws.Range("E22").Value = "  -2.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.49"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.31%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("E28").Value = "  -4.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.95%  "
$ws.Range("E30").Value = "  -0.64%  "
$ws.Range("E31").Value = "  -2.50%  "
$ws.Range("E32").Value = "  -3.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.670"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +23.95%  "
$ws.Range("E34").Value = "  -1.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.313.87"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.51"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.27%  "
$ws.Range("E37").Value = "  -0.73%  "
$ws.Range("E38").Value = "  -1.09%  "
$ws.Range("E39").Value = "  -2.69%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.788"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.733.92"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E47").Value = "  +0.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.807"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.67%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0510"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.59%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0975"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.80%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.998"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.25%  "
